$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H48").Value = 2311.111
$ws.Range("J48").Value = 2311.111
$ws.Range("L48").Value = 6933.333
$ws.Range("N48").Value = -7517.333

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H56").Value = 2311.111
$ws.Range("J56").Value = 2311.111
$ws.Range("L56").Value = 6933.333
$ws.Range("N56").Value = -8001.333

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H98").Value = 2889.9092
$ws.Range("I98").Value = 2542.875
$ws.Range("J98").Value = 3815.3333
$ws.Range("K98").Value = 2542.875
$ws.Range("L98").Value = 3815.3333
$ws.Range("M98").Value = -1044.875
$ws.Range("N98").Value = -6811.3333

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H106").Value = 3529.5386
$ws.Range("I106").Value = 3908.6667
$ws.Range("K106").Value = 3908.6667
$ws.Range("M106").Value = -3277.6667

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H107").Value = 2383.182
$ws.Range("I107").Value = 2134.8333
$ws.Range("K107").Value = 2134.8333
$ws.Range("M107").Value = -214.8332999999998

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H122").Value = 2889.9092
$ws.Range("I122").Value = 2542.875
$ws.Range("J122").Value = 3815.3333
$ws.Range("K122").Value = 7628.625
$ws.Range("L122").Value = 11445.9999
$ws.Range("M122").Value = -5178.625
$ws.Range("N122").Value = -16345.9999

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H138").Value = 1348.49
$ws.Range("I138").Value = 730.875
$ws.Range("J138").Value = 1760.2333
$ws.Range("K138").Value = 2192.625
$ws.Range("L138").Value = 5280.699900000001
$ws.Range("M138").Value = 2947.375
$ws.Range("N138").Value = -15560.6999

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H141").Value = 750
$ws.Range("I141").Value = 750
$ws.Range("K141").Value = 2250
$ws.Range("M141").Value = 2930

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 17736.666
$ws.Range("I2").Value = 1235.6666
$ws.Range("K2").Value = 1235.6666
$ws.Range("M2").Value = -1122.6666

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H45").Value = 1076
$ws.Range("I45").Value = 1112.5454
$ws.Range("J45").Value = 875
$ws.Range("K45").Value = 1112.5454
$ws.Range("L45").Value = 875
$ws.Range("M45").Value = -735.5454
$ws.Range("N45").Value = -1629

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 1323.1666
$ws.Range("I61").Value = 987.9
$ws.Range("J61").Value = 2999.5
$ws.Range("K61").Value = 987.9
$ws.Range("L61").Value = 2999.5
$ws.Range("M61").Value = -775.9
$ws.Range("N61").Value = -3423.5

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H63").Value = 2181.7273
$ws.Range("J63").Value = 2033.3334
$ws.Range("L63").Value = 2033.3334
$ws.Range("N63").Value = -3405.3334

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H66").Value = 2181.7273
$ws.Range("J66").Value = 2033.3334
$ws.Range("L66").Value = 10166.667
$ws.Range("N66").Value = -17030.667

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H110").Value = 1481.9231
$ws.Range("I110").Value = 772.44446
$ws.Range("K110").Value = 772.44446
$ws.Range("M110").Value = 1272.55554

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H116").Value = 17736.666
$ws.Range("I116").Value = 1235.6666
$ws.Range("K116").Value = 1235.6666
$ws.Range("M116").Value = 1058.3334

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H122").Value = 1956.6
$ws.Range("I122").Value = 1951.7778
$ws.Range("K122").Value = 5855.3334
$ws.Range("M122").Value = -3405.3334

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H136").Value = 1323.1666
$ws.Range("I136").Value = 987.9
$ws.Range("J136").Value = 2999.5
$ws.Range("K136").Value = 2963.7
$ws.Range("L136").Value = 8998.5
$ws.Range("M136").Value = -413.6999999999998
$ws.Range("N136").Value = -14098.5

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 17736.666
$ws.Range("I3").Value = 1235.6666
$ws.Range("K3").Value = 1235.6666
$ws.Range("M3").Value = -1121.6666

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H99").Value = 55557244
$ws.Range("I99").Value = 71430130
$ws.Range("J99").Value = 2150
$ws.Range("K99").Value = 71430130
$ws.Range("L99").Value = 2150
$ws.Range("M99").Value = -71428632
$ws.Range("N99").Value = -5146

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H107").Value = 2324.8
$ws.Range("I107").Value = 1652.75
$ws.Range("J107").Value = 5013
$ws.Range("K107").Value = 1652.75
$ws.Range("L107").Value = 5013
$ws.Range("M107").Value = 267.25
$ws.Range("N107").Value = -8853

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H134").Value = 7606.3687
$ws.Range("I134").Value = 823
$ws.Range("K134").Value = 2469
$ws.Range("M134").Value = 66

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H58").Value = 1027.0541
$ws.Range("I58").Value = 817.34485
$ws.Range("J58").Value = 1787.25
$ws.Range("K58").Value = 817.34485
$ws.Range("L58").Value = 1787.25
$ws.Range("M58").Value = -614.34485
$ws.Range("N58").Value = -2193.25

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H107").Value = 588.8570999999999
$ws.Range("I107").Value = 461.72223
$ws.Range("K107").Value = 461.72223
$ws.Range("M107").Value = 1458.27777

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H134").Value = 27779852
$ws.Range("I134").Value = 37039036
$ws.Range("J134").Value = 2300
$ws.Range("K134").Value = 111117108
$ws.Range("L134").Value = 6900
$ws.Range("M134").Value = -111114573
$ws.Range("N134").Value = -11970

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H136").Value = 1027.0541
$ws.Range("I136").Value = 817.34485
$ws.Range("J136").Value = 1787.25
$ws.Range("K136").Value = 2452.03455
$ws.Range("L136").Value = 5361.75
$ws.Range("M136").Value = 97.96545000000015
$ws.Range("N136").Value = -10461.75

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H34").Value = 1617.6154
$ws.Range("J34").Value = 2375
$ws.Range("L34").Value = 7125
$ws.Range("N34").Value = -7293

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H39").Value = 1592.6428
$ws.Range("J39").Value = 1684.1538
$ws.Range("L39").Value = 5052.4614
$ws.Range("N39").Value = -5640.4614

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H55").Value = 3374.75
$ws.Range("J55").Value = 3374.75
$ws.Range("L55").Value = 10124.25
$ws.Range("N55").Value = -10478.25

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H107").Value = 755.2222
$ws.Range("I107").Value = 0
$ws.Range("J107").Value = 755.2222
$ws.Range("K107").Value = 0
$ws.Range("L107").Value = 2265.6666
$ws.Range("M107").ClearContents()
$ws.Range("N107").Value = -6105.6666

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H108").Value = 871
$ws.Range("I108").Value = 369.5
$ws.Range("K108").Value = 1108.5
$ws.Range("M108").Value = 1771.5

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H102").Value = 2270.2354
$ws.Range("I102").Value = 2440.2
$ws.Range("J102").Value = 995.5
$ws.Range("K102").Value = 2440.2
$ws.Range("L102").Value = 995.5
$ws.Range("M102").Value = -818.1999999999998
$ws.Range("N102").Value = -4239.5

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H132").Value = 3028.7058
$ws.Range("I132").Value = 2917.5715
$ws.Range("K132").Value = 8752.7145
$ws.Range("M132").Value = -6222.7145

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H55").Value = 142.35715
$ws.Range("I55").Value = 58
$ws.Range("K55").Value = 58
$ws.Range("M55").Value = 115

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H61").Value = 1500
$ws.Range("I61").Value = 1000
$ws.Range("J61").Value = 2000
$ws.Range("K61").Value = 1000
$ws.Range("L61").Value = 2000
$ws.Range("M61").Value = -798
$ws.Range("N61").Value = -2404

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H113").Value = 1500
$ws.Range("I113").Value = 1000
$ws.Range("J113").Value = 2000
$ws.Range("K113").Value = 1000
$ws.Range("L113").Value = 2000
$ws.Range("M113").Value = 1170
$ws.Range("N113").Value = -6340

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H122").Value = 31253324
$ws.Range("I122").Value = 41670100
$ws.Range("J122").Value = 2997.5
$ws.Range("K122").Value = 125010300
$ws.Range("L122").Value = 8992.5
$ws.Range("M122").Value = -125007850
$ws.Range("N122").Value = -13892.5

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H110").Value = 20000
$ws.Range("J110").Value = 20000
$ws.Range("L110").Value = 20000
$ws.Range("N110").Value = -28180

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H113").Value = 649.3333
$ws.Range("I113").Value = 379.2
$ws.Range("J113").Value = 2000
$ws.Range("K113").Value = 1137.6
$ws.Range("L113").Value = 6000
$ws.Range("M113").Value = 1032.4
$ws.Range("N113").Value = -10340

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H136").Value = 827.65216
$ws.Range("I136").Value = 681.0526
$ws.Range("K136").Value = 2043.1578
$ws.Range("M136").Value = 506.8422
